$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Insert a new row at position 2, shifting existing rows down
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with data
$ws.Cells.Item(2, 1).Value = 117
$ws.Cells.Item(2, 2).Value = "https://leetcode.com/u/chestnut890123/"
$ws.Cells.Item(2, 4).Value = 2571
$ws.Cells.Item(2, 5).Value = 18
$ws.Cells.Item(2, 6).Value = 91
$ws.Cells.Item(2, 12).Value = "https://github.com/acmilannesta"

# Selection matches the target (L3)
$ws.Range("L3").Select()
